$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3586754328301822
$ws.Range("C2").Value = 0.06782755789456019
$ws.Range("E2").Value = 0.4157262113152598
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.5752500374835705
$ws.Range("H2").Value = 0.6916346986531892
$ws.Range("I2").Value = 0.529599250877034
$ws.Range("K2").Value = 0.3971614532771071
$ws.Range("N2").Value = 1.302844417335611
$ws.Range("B3").Value = 0.3182533808712549
$ws.Range("C3").Value = 0.05904574915909677
$ws.Range("E3").Value = 0.3627153676308268
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.5732819939856597
$ws.Range("H3").Value = 0.6954822043658453
$ws.Range("I3").Value = 0.5335095302643609
$ws.Range("K3").Value = 0.3496537670478119
$ws.Range("N3").Value = 1.320296940226546
$ws.Range("B4").Value = 0.2934875089234765
$ws.Range("C4").Value = 0.05364664019596432
$ws.Range("E4").Value = 0.3302567501748968
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.572556622140084
$ws.Range("H4").Value = 0.6982252203953863
$ws.Range("I4").Value = 0.5362553125640375
$ws.Range("K4").Value = 0.3205220898519485
$ws.Range("N4").Value = 1.331546672064095
$ws.Range("B5").Value = 0.283408797938506
$ws.Range("C5").Value = 0.05144458111044514
$ws.Range("E5").Value = 0.3170505433176913
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.5723819504186309
$ws.Range("H5").Value = 0.699438589171038
$ws.Range("I5").Value = 0.537460745669442
$ws.Range("K5").Value = 0.3086603049396217
$ws.Range("N5").Value = 1.336265115432298
$ws.Range("B6").Value = 0.2817360604859402
$ws.Range("C6").Value = 0.051078814321869
$ws.Range("E6").Value = 0.3148588766810292
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.5723602341773812
$ws.Range("H6").Value = 0.6996458359506335
$ws.Range("I6").Value = 0.5376661253986406
$ws.Range("K6").Value = 0.3066912450605344
$ws.Range("N6").Value = 1.337056703655449
$ws.Range("B7").Value = 0.2933515287082002
$ws.Range("C7").Value = 0.05361695020110346
$ws.Range("E7").Value = 0.3300785641473709
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.5725537775538214
$ws.Range("H7").Value = 0.6982411975588008
$ws.Range("I7").Value = 0.5362712194869701
$ws.Range("K7").Value = 0.3203620787624857
$ws.Range("N7").Value = 1.331609763929323
$ws.Range("B8").Value = 0.3447268942850883
$ws.Range("C8").Value = 0.06480098458823136
$ws.Range("E8").Value = 0.3974282677437913
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.5744708814055315
$ws.Range("H8").Value = 0.6928822296734438
$ws.Range("I8").Value = 0.5308758122002182
$ws.Range("K8").Value = 0.3807728490457691
$ws.Range("N8").Value = 1.308751168574948
$ws.Range("B9").Value = 0.4458966574308647
$ws.Range("C9").Value = 0.08668286616023124
$ws.Range("E9").Value = 0.5303022766164673
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.5820881317705471
$ws.Range("H9").Value = 0.6854004767170636
$ws.Range("I9").Value = 0.5230415556843475
$ws.Range("K9").Value = 0.4995480359200428
$ws.Range("N9").Value = 1.26816689417067
$ws.Range("B10").Value = 0.5204902423663214
$ws.Range("C10").Value = 0.1027387290156696
$ws.Range("E10").Value = 0.6285405817227314
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.590073101568052
$ws.Range("H10").Value = 0.6817594636161175
$ws.Range("I10").Value = 0.518973889610038
$ws.Range("K10").Value = 0.5870190988513855
$ws.Range("N10").Value = 1.240942236859192
$ws.Range("B11").Value = 0.5544840651353127
$ws.Range("C11").Value = 0.1100406493594051
$ws.Range("E11").Value = 0.6733935968061928
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.5942323985227631
$ws.Range("H11").Value = 0.6805083726857646
$ws.Range("I11").Value = 0.5174930585981556
$ws.Range("K11").Value = 0.626861377670707
$ws.Range("N11").Value = 1.229121525549232
$ws.Range("B12").Value = 0.5673653740058739
$ws.Range("C12").Value = 0.1128055592120916
$ws.Range("E12").Value = 0.6904039529248678
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.5958838126823736
$ws.Range("H12").Value = 0.6800930738343709
$ws.Range("I12").Value = 0.516985701178136
$ws.Range("K12").Value = 0.6419561532974853
$ws.Range("N12").Value = 1.224726592713068
$ws.Range("B13").Value = 0.5645907743213456
$ws.Range("C13").Value = 0.1122100920747187
$ws.Range("E13").Value = 0.6867393037647531
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.595524744345326
$ws.Range("H13").Value = 0.6801799126780281
$ws.Range("I13").Value = 0.5170925905252943
$ws.Range("K13").Value = 0.6387048954917702
$ws.Range("N13").Value = 1.225669500073991
$ws.Range("B14").Value = 0.5555436470012864
$ws.Range("C14").Value = 0.1102681228809104
$ws.Range("E14").Value = 0.6747925270459092
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.5943667270981194
$ws.Range("H14").Value = 0.6804730330260611
$ws.Range("I14").Value = 0.51745024627823
$ws.Range("K14").Value = 0.6281030862403441
$ws.Range("N14").Value = 1.228758320834094
$ws.Range("B15").Value = 0.5500031337238624
$ws.Range("C15").Value = 0.1090785921093982
$ws.Range("E15").Value = 0.6674781586509653
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.5936673730840596
$ws.Range("H15").Value = 0.6806601969209538
$ws.Range("I15").Value = 0.5176762828315873
$ws.Range("K15").Value = 0.6216101317195069
$ws.Range("N15").Value = 1.230660908909655
$ws.Range("B16").Value = 0.5182698762519067
$ws.Range("C16").Value = 0.1022615014301493
$ws.Range("E16").Value = 0.6256127939310403
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.5898119364061216
$ws.Range("H16").Value = 0.6818493957996736
$ws.Range("I16").Value = 0.519078124056584
$ws.Range("K16").Value = 0.5844163433916947
$ws.Range("N16").Value = 1.241726114621937
$ws.Range("B17").Value = 0.4988180141434952
$ws.Range("C17").Value = 0.0980790130882383
$ws.Range("E17").Value = 0.5999730925466338
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.5875821352662456
$ws.Range("H17").Value = 0.6826828538213334
$ws.Range("I17").Value = 0.5200329466152454
$ws.Range("K17").Value = 0.5616122839893194
$ws.Range("N17").Value = 1.248658850815575
$ws.Range("B18").Value = 0.4876355522254983
$ws.Range("C18").Value = 0.09567316383606794
$ws.Range("E18").Value = 0.5852411270508782
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.5863491695539835
$ws.Range("H18").Value = 0.6832003668338587
$ws.Range("I18").Value = 0.5206169055144869
$ws.Range("K18").Value = 0.5485008161335259
$ws.Range("N18").Value = 1.252699442658697
$ws.Range("B19").Value = 0.4838503526652289
$ws.Range("C19").Value = 0.09485854806970906
$ws.Range("E19").Value = 0.5802557083833477
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.5859402014241795
$ws.Range("H19").Value = 0.6833821309287202
$ws.Range("I19").Value = 0.52082058797982
$ws.Range("K19").Value = 0.5440623270333447
$ws.Range("N19").Value = 1.25407662495645
$ws.Range("B20").Value = 0.5008881057323435
$ws.Range("C20").Value = 0.09852426521896973
$ws.Range("E20").Value = 0.60270088390034
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.5878143680239845
$ws.Range("H20").Value = 0.6825901830962096
$ws.Range("I20").Value = 0.51992770367535
$ws.Range("K20").Value = 0.5640393158045072
$ws.Range("N20").Value = 1.247915355541054
$ws.Range("B21").Value = 0.5582007758426357
$ws.Range("C21").Value = 0.1108385299687882
$ws.Range("E21").Value = 0.678300877678609
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.5947047870421187
$ws.Range("H21").Value = 0.6803853482909261
$ws.Range("I21").Value = 0.5173437427608434
$ws.Range("K21").Value = 0.6312168945891585
$ws.Range("N21").Value = 1.227848850665051
$ws.Range("B22").Value = 0.5957078925407302
$ws.Range("C22").Value = 0.1188856836587888
$ws.Range("E22").Value = 0.7278593420332555
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.5996535544105797
$ws.Range("H22").Value = 0.6792852139548415
$ws.Range("I22").Value = 0.5159663585764704
$ws.Range("K22").Value = 0.6751642831062838
$ws.Range("N22").Value = 1.215208296638609
$ws.Range("B23").Value = 0.5756851180108526
$ws.Range("C23").Value = 0.1145908130932867
$ws.Range("E23").Value = 0.7013947373610279
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.5969713394270713
$ws.Range("H23").Value = 0.6798411244279947
$ws.Range("I23").Value = 0.5166729172557112
$ws.Range("K23").Value = 0.6517048031015804
$ws.Range("N23").Value = 1.22191134625682
$ws.Range("B24").Value = 0.499952215218741
$ws.Range("C24").Value = 0.09832297073145924
$ws.Range("E24").Value = 0.601467622700369
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.5877092230742562
$ws.Range("H24").Value = 0.6826319600881163
$ws.Range("I24").Value = 0.5199751749400363
$ws.Range("K24").Value = 0.562942058183836
$ws.Range("N24").Value = 1.248251318840978
$ws.Range("B25").Value = 0.4184816427561486
$ws.Range("C25").Value = 0.08076758921784233
$ws.Range("E25").Value = 0.4942574579925036
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.5796106736476929
$ws.Range("H25").Value = 0.6870993453430856
$ws.Range("I25").Value = 0.5248654933430927
$ws.Range("K25").Value = 0.4673811827086922
$ws.Range("N25").Value = 1.278691247295995
